# Task #5: Implement end to end flow
#
# Updates the email addresses shown in the "validLogin" sheet:
#   B3: shimaa1@dxc.com -> shimaa2@dxc.com
#   B4: shimaa2@dxc.com -> shimaa3@dxc.com
# and leaves the active selection on B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("validLogin")

$ws.Range("B3").Value = "shimaa2@dxc.com"
$ws.Range("B4").Value = "shimaa3@dxc.com"

$ws.Activate()
$ws.Range("B4").Select()
